$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 151, shifting existing rows 151:268 down to 152:269
$ws.Rows("151").Insert()

# Populate the newly inserted row 151 with the new record's data
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value = "Ñuble"
$ws.Range("D151").Value = 44634
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 100114013
$ws.Range("G151").Value = "Zanahoria"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 120
$ws.Range("K151").Value = 6500
$ws.Range("L151").Value = 7000
$ws.Range("M151").Value = 6750
$ws.Range("N151").Value = "$/saco 20 kilos"
$ws.Range("O151").Value = "Provincia de Diguillín"
$ws.Range("P151").Value = 338
$ws.Range("Q151").Value = 20
$ws.Range("R151").Value = "Hortaliza"
